# Apply "Here Comes the Rain" edit:
# Insert a new policy-assignment row into the "ALZ Policy Assignments 03CY23"
# worksheet, right after row 16 (i.e. becomes the new row 17), pushing all
# subsequent rows down by one. This also grows the table/filter/dimension
# ranges from A1:J50 to A1:J51.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALZ Policy Assignments 03CY23")

# Insert a new row at position 17 by copying the formatting of row 16 (which
# already carries the wrap/alignment/number-format combination used by the
# new entry) and inserting that copy above the current row 17. Everything
# currently at/after row 17 shifts down to row 18 and beyond, and all the
# named ranges / autofilter / hyperlink anchors that Excel maintains
# automatically follow along.
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()
$excel.CutCopyMode = $false

# Populate the newly inserted row 17 with the new policy assignment entry.
$ws.Cells.Item(17, 1).Value = "Intermediate Root"
$ws.Cells.Item(17, 2).Value = "Resource Group and Resource locations should match"
$ws.Cells.Item(17, 3).Value = "Resource Group and Resource locations should match"
$ws.Cells.Item(17, 4).Value = "Policy"
$ws.Cells.Item(17, 5).Value = "Custom"
$ws.Cells.Item(17, 6).Value = "In order to improve resilience and reliability, you need to be aware of where resources are deployed. To aid this awareness, ensure that the location of the resource group matches the location of the resources it contains."
$ws.Cells.Item(17, 7).Value = "Audit, Deny"
$ws.Cells.Item(17, 8).Value = "Audit-ResourceRGLocation.json"
$ws.Cells.Item(17, 9).Value = "TBD"
$ws.Cells.Item(17, 10).Value = Get-Date -Year 2023 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0

# Row height for the newly inserted row.
$ws.Rows.Item(17).RowHeight = 72

# Grow the AutoFilter range so it covers the new row, the same way Excel
# extends a filtered range that had a row inserted inside of it.
$ws.AutoFilterMode = $false
$ws.Range("A1:J51").AutoFilter()

# The hidden _xlnm._FilterDatabase defined name backing the AutoFilter also
# needs to be grown to the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='ALZ Policy Assignments 03CY23'!`$A`$1:`$J`$51"
    }
}

# Update the view: scroll so row 13 is at the top, and select K17.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("K17").Select()
